$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Sunday)
$ws.Range("B2").Value = "CIE4818-Darwin Nicolas-Hall 2"
$ws.Range("C2").Value = "GEN0807-Osvaldo Boyle PhD-Hall 1`nGEN1809-Laisha Schultz-Hall 2"
$ws.Range("D2").Value = "GEN0801-Darwin Nicolas-Hall 2"
$ws.Range("E2").Value = "MEC0811-Dr. Yvonne Treutel PhD-Hall 1`nGEN0801-sec-Hall 2"
$ws.Range("F2").Value = "GEN1801-Prof. Eladio Franecki-Hall 1`nCIE1808-sec-Hall 2"

# Row 3 (Monday)
$ws.Range("B3").Value = "POW1804-sec-Hall 1`nGEN0807-sec-Hall 2"
$ws.Range("C3").Value = "CIE1803-Vernie Sporer-Hall 2"
$ws.Range("D3").Value = "CIE2802-Destinee Feest-Hall 1"
$ws.Range("E3").Value = "CIE4818-sec-Hall 2"
$ws.Range("F3").Value = "CIE1803-lab-Hall 1"

# Row 4 (Tuesday)
$ws.Range("B4").Value = "GEN1801-Prof. Eladio Franecki-Hall 1`nGEN0802-lab-Hall 2"
$ws.Range("C4").Value = "GEN0809-Darwin Nicolas-Hall 1"
$ws.Range("D4").Value = "POW1804-Dr. Yvonne Treutel PhD-Hall 2"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "POW1804-Dr. Yvonne Treutel PhD-Hall 1"

# Row 5 (Wednesday)
$ws.Range("B5").Value = "GEN0810-Dr. Yvonne Treutel PhD-Hall 1"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "GEN1801-sec-Hall 1"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "MEC0811-sec-Hall 2"

# Row 6 (Thursday)
$ws.Range("B6").Value = "CIE3804-Vernie Sporer-Hall 1`nGEN0806-sec-Hall 2"
$ws.Range("C6").Value = "GEN1805-Mr. Howard Willms II-Hall 1"
$ws.Range("D6").Value = "GEN0801-Mr. Howard Willms II-Hall 2"
$ws.Range("E6").Value = "GEN0802-Osvaldo Boyle PhD-Hall 2"
$ws.Range("F6").Value = "CIE2802-sec-Hall 1"

# Row 7 (Saturday)
$ws.Range("B7").Value = "GEN2810-Mr. Howard Willms II-Hall 2"
$ws.Range("C7").Value = "GEN1809-sec-Hall 1"
$ws.Range("D7").Value = "CIE1808-Destinee Feest-Hall 2"
$ws.Range("E7").Value = "CIE3804-lab-Hall 1`nGEN0810-Dr. Yvonne Treutel PhD-Hall 2"
$ws.Range("F7").Value = "CIE1808-Destinee Feest-Hall 1`nGEN0806-Osvaldo Boyle PhD-Hall 2"
